$d = $word.ActiveDocument

# 1. Insert a "_GoBack" bookmark between ". Then" and " draw your traces." in the
#    Introduction paragraph.
$rBm = $d.Content
$rBm.Find.Execute(". Then", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rBm.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rBm) | Out-Null

# 2. Remove the degree-suffix runs from the "Page 1/2/3" table headers.
$d.Content.Find.Execute(" (180°)", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute(" (90°)", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute(" (45°)", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 3. Grow the "Speed & direction" row's height.
$tbl = $d.Tables.Item(1)
$speedRow = $tbl.Rows.Item(2)
$speedRow.Height = 98.95

# 4. Remove the "Traces" runs from the 2nd and 3rd cells of the Traces row,
#    leaving the 1st cell's "Traces" run untouched. Use direct Range deletion
#    (rather than Find) so duplicate "Traces" text in sibling cells isn't
#    accidentally matched.
$tracesCell2 = $tbl.Cell(3, 2)
$d.Range($tracesCell2.Range.Start, $tracesCell2.Range.End - 1).Delete()

$tracesCell3 = $tbl.Cell(3, 3)
$d.Range($tracesCell3.Range.Start, $tracesCell3.Range.End - 1).Delete()

# 5. Remove the "As you did this activity..." paragraph along with the two
#    blank paragraphs and the (now redundant) bookmark paragraph that used to
#    follow it, leaving the single blank paragraph right before the next
#    question ("If you set ...") intact. (Note: the Document.Paragraphs
#    collection can't be used reliably once a Table has been touched in this
#    runtime, so locate everything with Find + absolute Range offsets
#    instead.)
$rStart = $d.Content
$rStart.Find.Execute("As you did this activity", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deleteFrom = $rStart.Start

$rEnd = $d.Content
$rEnd.Find.Execute("If you set", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deleteTo = $rEnd.Start - 1

$d.Range($deleteFrom, $deleteTo).Delete()
